$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix rows 318-319: move the "Mode of Payment" value from column I to column H ---
$i318 = $ws.Range("I318").Value2
$ws.Range("H318").Value = $i318
$ws.Range("I318").ClearContents()

$i319 = $ws.Range("I319").Value2
$ws.Range("H319").Value = $i319
$ws.Range("I319").ClearContents()

# --- Fix row 320: phone number becomes a real number (was stored as text), and move Mode of Payment from I to H ---
$ws.Range("C320").Value = 9000171648
$i320 = $ws.Range("I320").Value2
$ws.Range("H320").Value = $i320
$ws.Range("I320").ClearContents()

# --- New date-separator row 321 (matches the style used for the date row at A300) ---
$ws.Range("A321").Value = 45843
$ws.Range("A321").NumberFormat = $ws.Range("A300").NumberFormat

# --- New fresher records: rows 322-336 ---
$ws.Range("A322").Value = 'YEDAMA RAKESH'
$ws.Range("B322").Value = '24C71A0557'
$ws.Range("C322").Value = 8309143154
$ws.Range("D322").Value = 'rakeshvarma5802@gmail.com'
$ws.Range("E322").Value = 'CSE'
$ws.Range("F322").Value = '1st Year'
$ws.Range("G322").Value = 'DRKFEODENC09'
$ws.Range("H322").Value = 'Cash'
$ws.Range("A323").Value = 'SHAIK SAMEERA BEGUM'
$ws.Range("B323").Value = '23C71A0408'
$ws.Range("C323").Value = 9014911140
$ws.Range("D323").Value = 'sameerashaik23c7@gmail.com'
$ws.Range("E323").Value = 'ECE'
$ws.Range("F323").Value = '2nd Year'
$ws.Range("G323").Value = 'OEL6ABUKIBYH'
$ws.Range("H323").Value = 'UPI'
$ws.Range("A324").Value = 'TELUKULA MOUNIKA'
$ws.Range("B324").Value = '23C71A0426'
$ws.Range("C324").Value = 8919285911
$ws.Range("D324").Value = 'mounikasahu878@gmail.com'
$ws.Range("E324").Value = 'ECE'
$ws.Range("F324").Value = '2nd Year'
$ws.Range("G324").Value = '4RPI38ROZ5EQ'
$ws.Range("H324").Value = 'UPI'
$ws.Range("A325").Value = 'VELLAPU GANESH'
$ws.Range("B325").Value = '24C71A6726'
$ws.Range("C325").Value = 9985969919
$ws.Range("D325").Value = 'kv.ganesh998@gmail.com'
$ws.Range("E325").Value = 'DS'
$ws.Range("F325").Value = '1st Year'
$ws.Range("G325").Value = 'IFX55G39KIK2'
$ws.Range("H325").Value = 'Cash'
$ws.Range("A326").Value = 'BEKKARIPALLI PAVITHRA'
$ws.Range("B326").Value = '24C71A6258'
$ws.Range("C326").Value = 6281492501
$ws.Range("D326").Value = 'pavithrabekkaripally26@gmail.com'
$ws.Range("E326").Value = 'CS'
$ws.Range("F326").Value = '1st Year'
$ws.Range("G326").Value = '8AT98EJM5H42'
$ws.Range("H326").Value = 'UPI'
$ws.Range("A327").Value = 'BOKKA RITHVIK'
$ws.Range("B327").Value = '24C71A66F1'
$ws.Range("C327").Value = 7207419210
$ws.Range("D327").Value = 'rithvikrishi73@gmail.com'
$ws.Range("E327").Value = 'CSE'
$ws.Range("F327").Value = '1st Year'
$ws.Range("G327").Value = '30AXX61JJIOJ'
$ws.Range("H327").Value = 'UPI'
$ws.Range("A328").Value = 'YERRA PAVAN KUMAR'
$ws.Range("B328").Value = '24C71A6211'
$ws.Range("C328").Value = 7780565324
$ws.Range("D328").Value = 'yerrapavanyadav@gmail.com'
$ws.Range("E328").Value = 'CS'
$ws.Range("F328").Value = '1st Year'
$ws.Range("G328").Value = '5INP0NKMLNGD'
$ws.Range("H328").Value = 'UPI'
$ws.Range("A329").Value = 'PARLAGOLLA SRIDHAR'
$ws.Range("B329").Value = '24C71A05B4'
$ws.Range("C329").Value = 9550598958
$ws.Range("D329").Value = 'r15omeo143@gmail.com'
$ws.Range("E329").Value = 'CSE'
$ws.Range("F329").Value = '1st Year'
$ws.Range("G329").Value = '38XA2YAGBN18'
$ws.Range("H329").Value = 'UPI'
$ws.Range("A330").Value = 'KUMMARI PRASHANTH'
$ws.Range("B330").Value = '24C71A0597'
$ws.Range("C330").Value = 8074918042
$ws.Range("D330").Value = 'Pkummari74@gmail.com'
$ws.Range("E330").Value = 'CSE'
$ws.Range("F330").Value = '1st Year'
$ws.Range("G330").Value = 'XAH2NOXND65P'
$ws.Range("H330").Value = 'UPI'
$ws.Range("A331").Value = 'SEETHA BALU CHANDER YADAV'
$ws.Range("B331").Value = '23C71A6671'
$ws.Range("C331").Value = 9701178486
$ws.Range("D331").Value = 'seethabaluchanderyadav5122@gmail.com'
$ws.Range("E331").Value = 'AI&ML'
$ws.Range("F331").Value = '2nd Year'
$ws.Range("G331").Value = '3XJK2G8SWBTB'
$ws.Range("H331").Value = 'UPI'
$ws.Range("A332").Value = 'THORLIKONDA VIVEK'
$ws.Range("B332").Value = '24C71A0303'
$ws.Range("C332").Value = 8885714218
$ws.Range("D332").Value = 'tvivekgaming96@gmail.com'
$ws.Range("E332").Value = 'MECH'
$ws.Range("F332").Value = '1st Year'
$ws.Range("G332").Value = '21JZR0X0TPUR'
$ws.Range("H332").Value = 'UPI'
$ws.Range("A333").Value = 'M Revanth sai'
$ws.Range("B333").Value = '24C71A0403'
$ws.Range("C333").Value = 6305383718
$ws.Range("D333").Value = 'revanthsai604@gmail.com'
$ws.Range("E333").Value = 'ECE'
$ws.Range("F333").Value = '1st Year'
$ws.Range("G333").Value = '8S7XN091S1XA'
$ws.Range("H333").Value = 'UPI'
$ws.Range("A334").Value = 'GAJULA AKHIL RAJU'
$ws.Range("B334").Value = '24C71A0424'
$ws.Range("C334").Value = 6305472071
$ws.Range("D334").Value = 'akhilrajgajula@gmail.com'
$ws.Range("E334").Value = 'ECE'
$ws.Range("F334").Value = '1st Year'
$ws.Range("G334").Value = 'C866ZYZV6SS7'
$ws.Range("H334").Value = 'Cash'
$ws.Range("A335").Value = 'D anji'
$ws.Range("B335").Value = '24C71A0408'
$ws.Range("C335").Value = 8074277041
$ws.Range("D335").Value = 'anji.dadapuram77@gmail.com'
$ws.Range("E335").Value = 'ECE'
$ws.Range("F335").Value = '1st Year'
$ws.Range("G335").Value = 'IB7Y11V6TB0J'
$ws.Range("H335").Value = 'UPI'
$ws.Range("A336").Value = 'GUDA GOUTHAMI'
$ws.Range("B336").Value = '24C71A05E4'
$ws.Range("C336").Value = 7989021353
$ws.Range("D336").Value = 'gouthamiguda2007@gmail.com'
$ws.Range("E336").Value = 'CSE'
$ws.Range("F336").Value = '1st Year'
$ws.Range("G336").Value = 'C5ZX7S6STDQR'
$ws.Range("I336").Value = 'Cash'
# --- Column A width ---
$ws.Columns("A").ColumnWidth = 8

# --- Final selection, matching the saved view state ---
$ws.Range("K332").Select()
